$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.494.43"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.741.18"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.87"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4509"
$ws.Range("E7").Value = "  +6.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3523"
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07379"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.21"
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.076"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.40"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.897"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.055"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "1.742.57"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.51"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001053"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06354"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.76"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.717"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "27.543.73"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.099"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.53"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.04"
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("D28").Value = "1.940.97"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "124.80"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.027"
$ws.Range("E30").Value = "  -4.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.045"
$ws.Range("E31").Value = "  -5.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09058"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.643"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.369"
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02266"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("E36").Value = "  -4.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05986"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2057"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6227"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.882"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.182"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.376"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.689"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.18"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.702"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5783"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.91"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.924"
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06840"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("E50").Value = "  -5.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.12"
$ws.Range("E51").Value = "  -2.72%  "
